$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A6").Value = "Tanja"
$ws.Range("B6").Value = "tanja@fvv.ccc"

$ws.Range("A7").Value = "Olga"
$ws.Range("B7").Value = "olga@gcm.com"

$ws.Range("A8").Value = "Olga"
$ws.Range("B8").Value = "olga@tan.ccc"
